$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.697.32'
$ws.Range("E2").Value = '  -0.57%  '

$ws.Range("D3").Value = '1.867.17'
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7302'
$ws.Range("E5").Value = '  -0.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.77'
$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9989'
$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3129'
$ws.Range("E8").Value = '  -0.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07086'
$ws.Range("E9").Value = '  -0.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.38'
$ws.Range("E10").Value = '  -1.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08230'
$ws.Range("E11").Value = '  -2.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7456'
$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.845.37'
$ws.Range("E13").Value = '  -0.88%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.310'
$ws.Range("E14").Value = '  -1.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.26'
$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("D16").Value = '29.702.53'
$ws.Range("E16").Value = '  -0.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.019'
$ws.Range("E17").Value = '  -0.35%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '248.28'
$ws.Range("E18").Value = '  +2.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.34'
$ws.Range("E19").Value = '  -1.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007802'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9974'
$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.127.60'
$ws.Range("E22").Value = '  +0.79%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("E24").Value = '  -2.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1534'
$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.158'
$ws.Range("E26").Value = '  -1.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.70'
$ws.Range("E27").Value = '  -0.73%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.53'
$ws.Range("E28").Value = '  -0.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.013'
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.432'
$ws.Range("E30").Value = '  -2.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.537'
$ws.Range("E31").Value = '  -2.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.524'
$ws.Range("E32").Value = '  -0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.191'
$ws.Range("E33").Value = '  -1.97%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05269'
$ws.Range("E34").Value = '  -1.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.232'
$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7530'
$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9974'
$ws.Range("E37").Value = '  -0.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.692'
$ws.Range("E38").Value = '  -0.37%  '

$ws.Range("E39").Value = '  -1.03%  '

$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4470'
$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.992'
$ws.Range("E42").Value = '  -1.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8650'
$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.16'
$ws.Range("E44").Value = '  -1.46%  '

$ws.Range("D45").Value = '1.052.15'
$ws.Range("E45").Value = '  -4.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.82'
$ws.Range("E46").Value = '  +0.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9988'
$ws.Range("E47").Value = '  -0.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.821'
$ws.Range("E48").Value = '  -1.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.471'
$ws.Range("E49").Value = '  -3.11%  '

$ws.Range("D50").Value = '2.020.07'
$ws.Range("E50").Value = '  +0.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.884'
$ws.Range("E51").Value = '  -5.84%  '
